$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34/35: Aptos and Monero swap places (identity + metrics)
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"

$ws.Range("D2").Value = "60.853.29"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "3.367.32"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.70"
$ws.Range("D9").ClearFormats()
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.381"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("D12").Value = "3.941.60"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.09"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "3.349.36"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "60.957.30"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.50"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.90"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.55"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -5.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.188"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -4.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.92"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "167.16"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.81"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.91"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").Value = "3.404.11"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.60"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.73%  "
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "2.440.56"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.63"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.04"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.37%  "
$ws.Range("E49").Value = "  -5.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.96"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.21%  "
$ws.Range("E51").Value = "  -3.11%  "
